# Applies two edits to "Kitegirl abilities" per the commit:
#   1. Title "Kitegirl abilities" -> split so "Kitegirl" is wrapped in a
#      proofErr (spellStart/spellEnd) pair, as Word's spellchecker would
#      flag the made-up word, and " abilities" becomes its own run.
#   2. "...dealing (33%) of the attack's damage..." -> the flat 33% crit
#      number is replaced by an AD-scaling value "10% AD" (split into its
#      own run) so the sticky bomb no longer does a fixed chunk of damage
#      (i.e. "enemies don't stick to player anymore" - the bomb damage is
#      now scaled off attack damage instead of the old flat rate).
#
# Both edits are applied by rebuilding the affected paragraph's content via
# Range.InsertXML so the exact run layout from the target XML (including
# the bare, rsid-less <w:r> runs and the <w:proofErr/> markers) is
# reproduced, rather than relying on Find/Replace (which would just fold
# the new text into the paragraph's single pre-existing run).

$d = $word.ActiveDocument

function New-OoxmlPackage([string]$BodyXml) {
    # Wraps a <w:body> fragment in the pkg:package envelope Range.InsertXML
    # expects, scoped to the main document part.
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $BodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

function Set-ParagraphXml($Paragraph, [string]$NewParagraphInnerXml) {
    # Replaces a paragraph's run content in place, keeping its own <w:p>
    # attributes (paraId, rsids, ...) untouched. The range must exclude the
    # trailing paragraph mark (End - 1) or InsertXML anchors the new
    # content at the wrong spot.
    $start = $Paragraph.Range.Start
    $end = $Paragraph.Range.End
    $target = $d.Range($start, $end - 1)
    $xml = New-OoxmlPackage("<w:p>" + $NewParagraphInnerXml + "</w:p>")
    $target.InsertXML($xml)
}

# --- Edit 1: title paragraph ------------------------------------------------
$titleXml = '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Kitegirl</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> abilities</w:t></w:r>'

$titlePara = $d.Paragraphs.Item(1)
if ($titlePara.Range.Text.TrimEnd("`r") -eq "Kitegirl abilities") {
    Set-ParagraphXml $titlePara $titleXml
}

# --- Edit 2: sticky bomb passive paragraph ----------------------------------
$stickyXml = '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">[Passive] </w:t></w:r>' +
    '<w:r w:rsidR="00381958"><w:rPr><w:lang w:val="en-US"/></w:rPr>' +
    '<w:t>All basic attacks apply a sticky bomb that explodes after (3s) dealing (</w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>10% AD</w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr>' +
    '<w:t>) of the attack' + [char]0x2019 + 's damage in a (small) area</w:t></w:r>' +
    '<w:r w:rsidR="004D2B0F"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>.</w:t></w:r>'

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*sticky bomb*" -and $p.Range.Text -like "*(33%)*") {
        Set-ParagraphXml $p $stickyXml
        break
    }
}
